# Generate Report for Handoff
#
# The localization report is regenerated: the row describing
# "9e3ed300-8a29-46c9-a992-8d928d43bccb.md" now sorts before the row
# describing "2a79a01b-9591-4291-9c83-3c4b86e07296.md" (on every sheet),
# and the 2a79a01b file's handoff status/date/error detail are refreshed
# because it is "Ready for handoff" again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Capture the two existing hyperlink target addresses (order: B2 then B3)
# before we touch anything, so we can recreate them after swapping rows.
$ovAddrs = @()
foreach ($h in $wsOverview.Hyperlinks) {
    $ovAddrs += $h.Address()
}

$wsOverview.Range("A2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.md"
$wsOverview.Range("B2").Value = "e2e\9e3ed300-8a29-46c9-a992-8d928d43bccb.md"
$wsOverview.Range("A3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.md"
$wsOverview.Range("B3").Value = "e2e\2a79a01b-9591-4291-9c83-3c4b86e07296.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 09:37:39"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ovAddrs[1], "", "", "e2e\9e3ed300-8a29-46c9-a992-8d928d43bccb.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ovAddrs[0], "", "", "e2e\2a79a01b-9591-4291-9c83-3c4b86e07296.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Existing hyperlink order: A2, I2, A3, I3
$zhAddrs = @()
foreach ($h in $wsZh.Hyperlinks) {
    $zhAddrs += $h.Address()
}

$wsZh.Range("A2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.md"
$wsZh.Range("G2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.72d4b330eb3170799a51c827cdbbf03bca5ccd28.zh-cn.xlf"
$wsZh.Range("I2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.md"
$wsZh.Range("J2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.72d4b330eb3170799a51c827cdbbf03bca5ccd28.zh-cn.xlf"

$wsZh.Range("A3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.40a0a378bb2d066a83a19f74cc447d8e0e2b4f4c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-06 09:37:28"
$wsZh.Range("I3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.md"
$wsZh.Range("J3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.40a0a378bb2d066a83a19f74cc447d8e0e2b4f4c.zh-cn.xlf"
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5dd95cc1f8d839634af2788af23ac253e38c7972/e2e/2a79a01b-9591-4291-9c83-3c4b86e07296.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f2038daccdc52f17863d82dbbfdaf845746c1de/e2e/2a79a01b-9591-4291-9c83-3c4b86e07296.md."

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhAddrs[2], "", "", "9e3ed300-8a29-46c9-a992-8d928d43bccb.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhAddrs[3], "", "", "9e3ed300-8a29-46c9-a992-8d928d43bccb.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhAddrs[0], "", "", "2a79a01b-9591-4291-9c83-3c4b86e07296.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhAddrs[1], "", "", "2a79a01b-9591-4291-9c83-3c4b86e07296.md") | Out-Null

# Column P (Error Detail) widened to fit the new long error message.
$wsZh.Columns("P").ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Existing hyperlink order: A2, I2, A3, I3
$deAddrs = @()
foreach ($h in $wsDe.Hyperlinks) {
    $deAddrs += $h.Address()
}

$wsDe.Range("A2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.md"
$wsDe.Range("G2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.72d4b330eb3170799a51c827cdbbf03bca5ccd28.de-de.xlf"
$wsDe.Range("I2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.md"
$wsDe.Range("J2").Value = "9e3ed300-8a29-46c9-a992-8d928d43bccb.72d4b330eb3170799a51c827cdbbf03bca5ccd28.de-de.xlf"

$wsDe.Range("A3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.40a0a378bb2d066a83a19f74cc447d8e0e2b4f4c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-06 09:37:39"
$wsDe.Range("I3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.md"
$wsDe.Range("J3").Value = "2a79a01b-9591-4291-9c83-3c4b86e07296.40a0a378bb2d066a83a19f74cc447d8e0e2b4f4c.de-de.xlf"
$wsDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5dd95cc1f8d839634af2788af23ac253e38c7972/e2e/2a79a01b-9591-4291-9c83-3c4b86e07296.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f2038daccdc52f17863d82dbbfdaf845746c1de/e2e/2a79a01b-9591-4291-9c83-3c4b86e07296.md."

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deAddrs[2], "", "", "9e3ed300-8a29-46c9-a992-8d928d43bccb.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deAddrs[3], "", "", "9e3ed300-8a29-46c9-a992-8d928d43bccb.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deAddrs[0], "", "", "2a79a01b-9591-4291-9c83-3c4b86e07296.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deAddrs[1], "", "", "2a79a01b-9591-4291-9c83-3c4b86e07296.md") | Out-Null

# Column P (Error Detail) widened to fit the new long error message.
$wsDe.Columns("P").ColumnWidth = 39.14

Write-Output "Done"
